# Cascades the historical "Betarraga" price rows down by one week and
# appends the previous last row as a new row 167 (weekly fruit/veggie refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 152-166: cascade the date/volume/price/origin data down one row ---
$ws.Cells.Item(152, 4).Value = 44449

$ws.Cells.Item(153, 4).Value = 44161
$ws.Cells.Item(153, 11).Value = 600
$ws.Cells.Item(153, 12).Value = 600
$ws.Cells.Item(153, 13).Value = 600
$ws.Cells.Item(153, 15).Value = 'Región del Maule'
$ws.Cells.Item(153, 16).Value = 120

$ws.Cells.Item(154, 11).Value = 800
$ws.Cells.Item(154, 12).Value = 800
$ws.Cells.Item(154, 13).Value = 800
$ws.Cells.Item(154, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(154, 16).Value = 160

$ws.Cells.Item(155, 4).Value = 44438
$ws.Cells.Item(155, 10).Value = 3000
$ws.Cells.Item(155, 11).Value = 700
$ws.Cells.Item(155, 12).Value = 700
$ws.Cells.Item(155, 13).Value = 700
$ws.Cells.Item(155, 16).Value = 140

$ws.Cells.Item(156, 4).Value = 44251
$ws.Cells.Item(156, 10).Value = 2000

$ws.Cells.Item(157, 4).Value = 44428
$ws.Cells.Item(157, 10).Value = 3000

$ws.Cells.Item(158, 4).Value = 44435
$ws.Cells.Item(158, 10).Value = 12000

$ws.Cells.Item(159, 4).Value = 44376
$ws.Cells.Item(159, 10).Value = 2000
$ws.Cells.Item(159, 11).Value = 600
$ws.Cells.Item(159, 12).Value = 600
$ws.Cells.Item(159, 13).Value = 600
$ws.Cells.Item(159, 16).Value = 120

$ws.Cells.Item(160, 4).Value = 44279
$ws.Cells.Item(160, 11).Value = 500
$ws.Cells.Item(160, 12).Value = 500
$ws.Cells.Item(160, 13).Value = 500
$ws.Cells.Item(160, 16).Value = 100

$ws.Cells.Item(161, 4).Value = 44412

$ws.Cells.Item(162, 4).Value = 44223

$ws.Cells.Item(163, 4).Value = 44314
$ws.Cells.Item(163, 11).Value = 600
$ws.Cells.Item(163, 12).Value = 600
$ws.Cells.Item(163, 13).Value = 600
$ws.Cells.Item(163, 16).Value = 120

$ws.Cells.Item(164, 4).Value = 44448
$ws.Cells.Item(164, 11).Value = 650
$ws.Cells.Item(164, 12).Value = 650
$ws.Cells.Item(164, 13).Value = 650
$ws.Cells.Item(164, 16).Value = 130

$ws.Cells.Item(165, 4).Value = 44167
$ws.Cells.Item(165, 10).Value = 3000
$ws.Cells.Item(165, 11).Value = 500
$ws.Cells.Item(165, 12).Value = 500
$ws.Cells.Item(165, 13).Value = 500
$ws.Cells.Item(165, 16).Value = 100

$ws.Cells.Item(166, 4).Value = 44399

# --- Add new row 167 with the data that used to be in row 166 ---
$ws.Cells.Item(167, 1).Value = 5
$ws.Cells.Item(167, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(167, 3).Value = 'Maule'
$ws.Cells.Item(167, 4).Value = 44400
$ws.Cells.Item(167, 5).Value = 7
$ws.Cells.Item(167, 6).Value = 100114014
$ws.Cells.Item(167, 7).Value = 'Betarraga'
$ws.Cells.Item(167, 8).Value = 'Sin especificar'
$ws.Cells.Item(167, 9).Value = 'Primera'
$ws.Cells.Item(167, 10).Value = 4000
$ws.Cells.Item(167, 11).Value = 650
$ws.Cells.Item(167, 12).Value = 650
$ws.Cells.Item(167, 13).Value = 650
$ws.Cells.Item(167, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(167, 15).Value = 'Región del Maule'
$ws.Cells.Item(167, 16).Value = 130
$ws.Cells.Item(167, 17).Value = 5
$ws.Cells.Item(167, 18).Value = 'Hortaliza'

# Match the date number format used by the other rows in column D
$ws.Cells.Item(167, 4).NumberFormat = $ws.Cells.Item(166, 4).NumberFormat

